$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers): swap H1/I1 so that "SamplePortion" comes before "Result" ---
$ws.Range("H1").Value = "SamplePortion"
$ws.Range("I1").Value = "Result"

# --- Row 2 (type/unit row): add unit annotation to the float columns (H2/I2) ---
$ws.Range("H2").Value = "#float,  unit:mg"
$ws.Range("I2").Value = "#float,  unit:mg"

# --- Row 3 (new description row, French labels matching each header in row 1) ---
$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"
$ws.Range("H3").Value = "#PriseEssai"
$ws.Range("I3").Value = "#Resultat"
$ws.Range("J3").Value = "#NuméroLotReactif"
